$wb = $excel.ActiveWorkbook

$wsWalMart  = $wb.Worksheets.Item("WalMart")
$wsCVS      = $wb.Worksheets.Item("CVS")
$wsPetSmart = $wb.Worksheets.Item("PetSmart")

# --- Header row: bold the header cells on all three sheets, and rename the
# "QtyAvailable" header to "Qty" ---
$wsWalMart.Range("A1:D1").Font.Bold = $true
$wsWalMart.Range("C1").Value = "Qty"

$wsCVS.Range("A1:D1").Font.Bold = $true
$wsCVS.Range("C1").Value = "Qty"

$wsPetSmart.Range("A1:D1").Font.Bold = $true
$wsPetSmart.Range("C1").Value = "Qty"

# A new (empty but styled/bold) cell E1 appears next to the WalMart header row
$wsWalMart.Range("E1").Font.Bold = $true

# --- WalMart data changes: Water now shows 0 available, and Chips is out of
# stock ("NA") instead of a numeric quantity ---
$wsWalMart.Range("C2").Value = 0
$wsWalMart.Range("C4").Value = "NA"

# --- Cosmetic/page setup touch-ups that mirror the author's session ---
$wsWalMart.PageSetup.Orientation = 1
$wsCVS.PageSetup.Orientation = 1

$wsPetSmart.Columns.Item(4).ColumnWidth = 14.25

# --- Selection / active-cell bookkeeping per sheet ---
$wsCVS.Activate()
$wsCVS.Range("D1:D1048576").Select() | Out-Null

$wsPetSmart.Activate()
$wsPetSmart.Range("D1:D1048576").Select() | Out-Null

$wsWalMart.Activate()
$wsWalMart.Range("C4").Select() | Out-Null
